$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" column (C) date for rows 2-4 to 46063 (2026-02-10)
$ws.Range("C2").Value = 46063
$ws.Range("C3").Value = 46063
$ws.Range("C4").Value = 46063

# Swap rows 3 and 4 for columns A (Beteckning), B (Datum) and G (Area (ha))
$ws.Range("A3").Value = "A 35536-2025"
$ws.Range("B3").Value = 45856
$ws.Range("G3").Value = 4

$ws.Range("A4").Value = "A 36713-2023"
$ws.Range("B4").Value = 45153
$ws.Range("G4").Value = 0.7
